$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT even though it looks
# like a number (matches the workbook's existing convention of numeric
# strings such as "3", "1.0", "24320.00" stored as text in columns D and G).
function Set-TextValue($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# ---------------------------------------------------------------------
# 1. Insert a new row at position 10 - this pushes the existing rows
#    10-20 down to 11-21 (matching the row-shift seen across the whole
#    bottom half of the sheet) and creates a brand new, empty row 10.
# ---------------------------------------------------------------------
$ws.Rows("10:10").Insert()

# ---------------------------------------------------------------------
# 2. Populate the new row 10 with the "Medium point" line item.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "P. point"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 84
Set-TextValue "D10" "3"
$ws.Range("E10").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F10").Value = 472
Set-TextValue "G10" "39648.00"
$ws.Range("H10").Value = 0

# ---------------------------------------------------------------------
# 3. Quantity updates on rows that are NOT shifted (row 8 and row 9),
#    including the recalculated "upto date amount" for row 9.
# ---------------------------------------------------------------------
$ws.Range("C8").Value = 50

$ws.Range("C9").Value = 95
Set-TextValue "G9" "24320.00"

# ---------------------------------------------------------------------
# 4. Quantity / amount updates on the rows that were shifted down by
#    the insert in step 1 (old row 10 -> new row 11, old row 11 -> new
#    row 12, etc). Only the cells whose values actually differ from the
#    shifted-down originals are touched here.
# ---------------------------------------------------------------------
$ws.Range("C11").Value = 88
Set-TextValue "G11" "58256.00"

$ws.Range("C12").Value = 56

$ws.Range("C13").Value = 28
Set-TextValue "G13" "3808.00"

$ws.Range("C14").Value = 32
Set-TextValue "G14" "736.00"

$ws.Range("C15").Value = 36

$ws.Range("C16").Value = 5

$ws.Range("C17").Value = 26

# ---------------------------------------------------------------------
# 5. Grand-total summary rows. Row 19 ("Grand Total Rs.") and row 21
#    ("NET PAYABLE AMOUNT Rs.") need their amounts refreshed; row 20
#    ("Tender Premium @ 0%") already carries the correct values after
#    the shift.
# ---------------------------------------------------------------------
Set-TextValue "G19" "126768.00"
Set-TextValue "H19" "126768.00"

Set-TextValue "G21" "126768.00"
Set-TextValue "H21" "126768.00"
